$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Result column (A) to "Pass" for rows 2-16
for ($r = 2; $r -le 16; $r++) {
    $ws.Cells.Item($r, 1).Value = "Pass"
}

# Update Date column (B) with new timestamps for rows 2-16
$dates = @(
    "Mon Jul 17 21:16:23 EDT 2023",
    "Mon Jul 17 21:17:06 EDT 2023",
    "Mon Jul 17 21:17:53 EDT 2023",
    "Mon Jul 17 21:18:34 EDT 2023",
    "Mon Jul 17 21:19:15 EDT 2023",
    "Mon Jul 17 21:19:55 EDT 2023",
    "Mon Jul 17 21:20:36 EDT 2023",
    "Mon Jul 17 21:21:22 EDT 2023",
    "Mon Jul 17 21:22:10 EDT 2023",
    "Mon Jul 17 21:22:57 EDT 2023",
    "Mon Jul 17 21:23:39 EDT 2023",
    "Mon Jul 17 21:24:26 EDT 2023",
    "Mon Jul 17 21:25:08 EDT 2023",
    "Mon Jul 17 21:25:52 EDT 2023",
    "Mon Jul 17 21:26:37 EDT 2023"
)

for ($i = 0; $i -lt $dates.Length; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 2).Value = $dates[$i]
}

# Remove rows 17-22 (old "Withholding Tax" test cases no longer needed)
$ws.Range("A17:A22").EntireRow.Delete()

# Update the active selection to reflect the author's final cursor position
$ws.Range("E27").Select()
